$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(10, "food_beef", "Beef"),
    @(11, "food_eggs", "Eggs"),
    @(12, "food_fruits", "Fruits"),
    @(13, "food_lamb", "Lamb"),
    @(14, "food_mayo", "Mayonnaise"),
    @(15, "food_milk", "Milk"),
    @(16, "food_pork", "Pork"),
    @(17, "food_poultry", "Poultry"),
    @(18, "food_rice_cooked", "Rice (Cooked)"),
    @(19, "food_shellfish", "Shellfish"),
    @(20, "food_vegetables_sliced", "Vegetables (Sliced)"),
    @(21, "food_bread", "Bread (Plain)"),
    @(22, "food_cannedGoods", "Canned Goods"),
    @(23, "food_flour", "Flour")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = $item[2]
}

# Row 25 (Rice (Raw)) was typed before row 24 (Potatoes) in the original
# edit session, and within row 25 the value (B) was entered before the
# key (A) -- this ordering is preserved in the shared-string table.
$ws.Cells.Item(25, 2).Value = "Rice (Raw)"
$ws.Cells.Item(25, 1).Value = "food_rice_raw"

$ws.Cells.Item(24, 1).Value = "food_potato"
$ws.Cells.Item(24, 2).Value = "Potatoes (Raw)"

$ws.Range("B24").Select()
